$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data block for rows 2-13 (cyclic shift of player blocks by one group)
$data = @(
    @("Baker Mayfield", "Group1", 3466.666666666667, 60.9, 14, 21.66666666666667, 170),
    @("Baker Mayfield", "Group2", 3328.666666666667, 65.3888888888889, 10.44444444444444, 25.22222222222222, 158.2222222222222),
    @("Baker Mayfield", "Difference", -138, 4.488888888888901, -3.555555555555555, 3.55555555555555, -11.7777777777778),
    @("Josh Allen", "Group1", 4013.333333333333, 63.76666666666667, 11.33333333333333, 31, 202.6666666666667),
    @("Josh Allen", "Group2", 4106.666666666667, 64.46666666666667, 12.66666666666667, 30.66666666666667, 195.3333333333333),
    @("Josh Allen", "Difference", 93.33333333333348, 0.6999999999999957, 1.333333333333332, -0.3333333333333321, -7.333333333333314),
    @("Kirk Cousins", "Group1", 4029.666666666667, 67.66666666666667, 8.666666666666666, 31.33333333333333, 189.6666666666667),
    @("Kirk Cousins", "Group2", 3462, 67.43333333333334, 11.66666666666667, 21.66666666666667, 168.3333333333333),
    @("Kirk Cousins", "Difference", -567.6666666666665, -0.2333333333333343, 3, -9.666666666666664, -21.33333333333331),
    @("Lamar Jackson", "Group1", 2922, 64.96666666666667, 9.333333333333334, 26, 144.6666666666667),
    @("Lamar Jackson", "Group2", 3364, 65.39999999999999, 6, 27.33333333333333, 156.6666666666667),
    @("Lamar Jackson", "Difference", 442, 0.4333333333333229, -3.333333333333334, 1.333333333333332, 12)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
